$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be stored as text so numeric-looking price strings
# (e.g. "1.00", "43.219.87") keep their exact original formatting instead
# of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.219.87"
$ws.Range("E2").Value = "  +1.67%  "
# Row 3
$ws.Range("D3").Value = "2.383.78"
$ws.Range("E3").Value = "  +6.94%  "
# Row 4
$ws.Range("E4").Value = "  -0.29%  "
# Row 5
$ws.Range("D5").Value = "324.23"
$ws.Range("E5").Value = "  +9.77%  "
# Row 6
$ws.Range("D6").Value = "105.99"
$ws.Range("E6").Value = "  -5.89%  "
# Row 7
$ws.Range("E7").Value = "  +3.44%  "
# Row 8
$ws.Range("E8").Value = "  -0.11%  "
# Row 9
$ws.Range("D9").Value = "0.653"
$ws.Range("E9").Value = "  +8.49%  "
# Row 10
$ws.Range("D10").Value = "41.68"
$ws.Range("E10").Value = "  -4.96%  "
# Row 11
$ws.Range("E11").Value = "  +2.51%  "
# Row 12
$ws.Range("D12").Value = "8.57"
$ws.Range("E12").Value = "  -0.95%  "
# Row 13
$ws.Range("E13").Value = "  -3.11%  "
# Row 14
$ws.Range("D14").Value = "17.25"
$ws.Range("E14").Value = "  +15.04%  "
# Row 16
$ws.Range("D16").Value = "2.744.27"
$ws.Range("E16").Value = "  +6.95%  "
# Row 17
$ws.Range("D17").Value = "2.389.47"
# Row 18
$ws.Range("D18").Value = "43.168.90"
$ws.Range("E18").Value = "  +1.48%  "
# Row 19
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  +7.07%  "
# Row 20
$ws.Range("D20").Value = "0.0000108"
$ws.Range("E20").Value = "  +2.85%  "
# Row 21
$ws.Range("D21").Value = "76.56"
$ws.Range("E21").Value = "  +4.08%  "
# Row 22
$ws.Range("D22").Value = "276.76"
$ws.Range("E22").Value = "  +17.44%  "
# Row 23
$ws.Range("E23").Value = "  +1.47%  "
# Row 24
$ws.Range("E24").Value = "  +0.46%  "
# Row 25
$ws.Range("D25").Value = "9.68"
$ws.Range("E25").Value = "  +8.81%  "
# Row 26
$ws.Range("D26").Value = "11.72"
$ws.Range("E26").Value = "  +2.30%  "
# Row 27
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.09%  "
# Row 28
$ws.Range("D28").Value = "22.92"
$ws.Range("E28").Value = "  +7.15%  "
# Row 29
$ws.Range("D29").Value = "176.37"
$ws.Range("E29").Value = "  +0.62%  "
# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "37.85"
$ws.Range("E30").Value = "  +0.92%  "
# Row 31
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "3.21"
$ws.Range("E31").Value = "  +2.33%  "
# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -2.32%  "
# Row 33
$ws.Range("D33").Value = "0.0928"
$ws.Range("E33").Value = "  +4.99%  "
# Row 34
$ws.Range("D34").Value = "5.91"
$ws.Range("E34").Value = "  +4.29%  "
# Row 35
$ws.Range("E35").Value = "  +5.40%  "
# Row 36
$ws.Range("D36").Value = "4.83"
$ws.Range("E36").Value = "  -3.77%  "
# Row 37
$ws.Range("E37").Value = "  -1.77%  "
# Row 38
$ws.Range("E38").Value = "  -2.47%  "
# Row 39
$ws.Range("E39").Value = "  +2.48%  "
# Row 40
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  +17.96%  "
# Row 41
$ws.Range("E41").Value = "  +20.87%  "
# Row 42
$ws.Range("E42").Value = "  +0.60%  "
# Row 43
$ws.Range("D43").Value = "124.42"
$ws.Range("E43").Value = "  +21.96%  "
# Row 44
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "95.33"
$ws.Range("E44").Value = "  +72.54%  "
# Row 45
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "69.14"
$ws.Range("E45").Value = "  -4.18%  "
# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.02%  "
# Row 47
$ws.Range("D47").Value = "12.48"
$ws.Range("E47").Value = "  +0.94%  "
# Row 48
$ws.Range("D48").Value = "9.51"
$ws.Range("E48").Value = "  +12.53%  "
# Row 49
$ws.Range("E49").Value = "  +4.09%  "
# Row 50
$ws.Range("E50").Value = "  +1.38%  "
# Row 51
$ws.Range("E51").Value = "  +5.54%  "
